$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Cxcl12", "Ackr3", "ECs", 3, 1, 81.05837566666666, 243.175127, 0.3545816884225585, 0.3545816884225585, 2, 0.6666666666666666, 10.577143, 31.731429, 0.2954534031561419, 0.2954534031561419, 857.3660307740536, 7716.294276966481, 0.1047623665412957, 0.1047623665412957),
    @("ECs", "Cxcl12", "Ackr3", "FAPs", 3, 1, 81.05837566666666, 243.175127, 0.3545816884225585, 0.3545816884225585, 3, 1, 16.17621166666667, 48.528635, 0.4518532827901404, 0.4518532827901404, 1311.217442140183, 11800.95697926164, 0.1602188999310038, 0.1602188999310038),
    @("ECs", "Cxcl12", "Ackr3", "sCs", 3, 1, 81.05837566666666, 243.175127, 0.3545816884225585, 0.3545816884225585, 3, 1, 9.046344666666666, 27.139034, 0.2526933140537176, 0.2526933140537176, 733.2820044008131, 6599.538039607317, 0.08960042195025904, 0.08960042195025904),
    @("FAPs", "Cxcl12", "Ackr3", "ECs", 3, 1, 123.018252, 369.054756, 0.5381309351710768, 0.5381309351710768, 2, 0.6666666666666666, 10.577143, 31.731429, 0.2954534031561419, 0.2954534031561419, 1301.181643014036, 11710.63478712632, 0.1589926161398918, 0.1589926161398918),
    @("FAPs", "Cxcl12", "Ackr3", "FAPs", 3, 1, 123.018252, 369.054756, 0.5381309351710768, 0.5381309351710768, 3, 1, 16.17621166666667, 48.528635, 0.4518532827901404, 0.4518532827901404, 1989.96928321534, 17909.72354893806, 0.2431562296279793, 0.2431562296279793),
    @("FAPs", "Cxcl12", "Ackr3", "sCs", 3, 1, 123.018252, 369.054756, 0.5381309351710768, 0.5381309351710768, 3, 1, 9.046344666666666, 27.139034, 0.2526933140537176, 0.2526933140537176, 1112.865507882856, 10015.7895709457, 0.1359820894032057, 0.1359820894032057),
    @("M2", "Cxcl12", "Ackr3", "ECs", 3, 1, 0.3624666666666667, 1.0874, 0.001585573873230423, 0.001585573873230423, 2, 0.6666666666666666, 10.577143, 31.731429, 0.2954534031561419, 0.2954534031561419, 3.833861766066666, 34.5047558946, 0.0004684631968013938, 0.0004684631968013938),
    @("M2", "Cxcl12", "Ackr3", "FAPs", 3, 1, 0.3624666666666667, 1.0874, 0.001585573873230423, 0.001585573873230423, 3, 1, 16.17621166666667, 48.528635, 0.4518532827901404, 0.4518532827901404, 5.863337522111111, 52.770037699, 0.0007164467597254447, 0.0007164467597254447),
    @("M2", "Cxcl12", "Ackr3", "sCs", 3, 1, 0.3624666666666667, 1.0874, 0.001585573873230423, 0.001585573873230423, 3, 1, 9.046344666666666, 27.139034, 0.2526933140537176, 0.2526933140537176, 3.278998396844444, 29.5109855716, 0.0004006639167035848, 0.0004006639167035848),
    @("sCs", "Cxcl12", "Ackr3", "ECs", 3, 1, 24.16373066666667, 72.491192, 0.1057018025331343, 0.1057018025331344, 2, 0.6666666666666666, 10.577143, 31.731429, 0.2954534031561419, 0.2954534031561419, 255.5832346748186, 2300.249112073368, 0.03122995727815304, 0.03122995727815305),
    @("sCs", "Cxcl12", "Ackr3", "FAPs", 3, 1, 24.16373066666667, 72.491192, 0.1057018025331343, 0.1057018025331344, 3, 1, 16.17621166666667, 48.528635, 0.4518532827901404, 0.4518532827901404, 390.8776219203244, 3517.89859728292, 0.04776170647143193, 0.04776170647143194),
    @("sCs", "Cxcl12", "Ackr3", "sCs", 3, 1, 24.16373066666667, 72.491192, 0.1057018025331343, 0.1057018025331344, 3, 1, 9.046344666666666, 27.139034, 0.2526933140537176, 0.2526933140537176, 218.5934360431698, 1967.340924388528, 0.02671013878354936, 0.02671013878354937)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
